# skompletowane czesci, wybrane rezystory, kondensator przy wejsciu
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Arkusz1
$ws2 = $wb.Worksheets.Item(2)   # Arkusz2

# ---------------------------------------------------------------
# Arkusz1 ("czesci" main BOM)
# ---------------------------------------------------------------

# Row5: KTIR0711S reflective sensor - now chosen
$ws1.Range("A5").Value = "KTIR0711S"
$ws1.Range("B5").Value = 10
$ws1.Range("C5").Value = 1.95
# D5 already has formula =B5*C5 from before, leave it as-is

# Row6: used to hold the "Razem:" total; now becomes "mocowanie silnika" (motor mount) line item
$ws1.Range("A6").Value = "mocowanie silnika"
$ws1.Range("B6").Value = 1
$ws1.Range("C6").Value = 17.9
$ws1.Range("D6").Formula = "=B6*C6"

# Row7: new "Razem:" (total) row, highlighted yellow
$ws1.Range("A7").Value = "Razem:"
$ws1.Range("D7").Formula = "=SUM(D3:D6)"
$ws1.Range("A7:D7").Interior.Color = 65535

# Row11: stray bold formatting on an otherwise empty cell C11
$ws1.Range("C11").Font.Bold = $true

# Bottom notes list shrinks from 5 lines (rows 12-16) to 3 (rows 13,14,16):
#   row13: "Nie wybrane" (was row12's text)
#   row14: "bezpiecznik"  (unchanged)
#   row16: "przelacznik on/off" (was row15's text)
# Rows 12 and 15 are cleared out entirely; "gniazdo miniusb" and
# "rezystory np. dla rgb" are no longer referenced anywhere.
$ws1.Range("A12").ClearContents()
$ws1.Range("A13").Value = "Nie wybrane"
$ws1.Range("A14").Value = "bezpiecznik"
$ws1.Range("A15").ClearContents()
$ws1.Range("A16").Value = "przelacznik on/off"

# ---------------------------------------------------------------
# Arkusz2 (detailed parts list)
# ---------------------------------------------------------------

# Row4: LED Biala quantity/price unchanged, just re-select item text (shared string reorder only)
$ws2.Range("A4").Value = "LED Biała"

# Row5: LED RGB - quantity reduced from 3 to 1
$ws2.Range("B5").Value = 1

# Row6: Kwarc 16MHz -> Kwarc 8MHz, price changed
$ws2.Range("A6").Value = "Kwarc 8MHz"
$ws2.Range("C6").Value = 1.8

# Row9: Zlacze FFC/FPC quantity 3 -> 2
$ws2.Range("B9").Value = 2

# Row11: IRML0030 quantity 5 -> 3
$ws2.Range("B11").Value = 3

# Row13: R47k quantity 1 -> 2
$ws2.Range("B13").Value = 2

# New rows 17-29: additional resistors / capacitors / regulators picked
$ws2.Range("A17").Value = "R1k"
$ws2.Range("B17").Value = 1
$ws2.Range("C17").Value = 0.3
$ws2.Range("D17").Formula = "=B17*C17"

$ws2.Range("A18").Value = "R150k"
$ws2.Range("B18").Value = 1
$ws2.Range("C18").Value = 0.3
$ws2.Range("D18").Formula = "=B18*C18"

$ws2.Range("A19").Value = "R13k"
$ws2.Range("B19").Value = 1
$ws2.Range("C19").Value = 0.3
$ws2.Range("D19").Formula = "=B19*C19"

$ws2.Range("A20").Value = "R22k"
$ws2.Range("B20").Value = 1
$ws2.Range("C20").Value = 0.3
$ws2.Range("D20").Formula = "=B20*C20"

$ws2.Range("A21").Value = "R62k"
$ws2.Range("B21").Value = 1
$ws2.Range("C21").Value = 0.3
$ws2.Range("D21").Formula = "=B21*C21"

$ws2.Range("A22").Value = "USB typ B smd mini"
$ws2.Range("B22").Value = 3
$ws2.Range("C22").Value = 1.01
$ws2.Range("D22").Formula = "=B22*C22"

$ws2.Range("A23").Value = "C22pF"
$ws2.Range("B23").Value = 1
$ws2.Range("C23").Value = 0.39
$ws2.Range("D23").Formula = "=B23*C23"

$ws2.Range("A24").Value = "C100nF"
$ws2.Range("B24").Value = 2
$ws2.Range("C24").Value = 0.39
$ws2.Range("D24").Formula = "=B24*C24"

$ws2.Range("A25").Value = "C10nF"
$ws2.Range("B25").Value = 1
$ws2.Range("C25").Value = 0.45
$ws2.Range("D25").Formula = "=B25*C25"

$ws2.Range("A26").Value = "C4,7"
$ws2.Range("B26").Value = 2
$ws2.Range("C26").Value = 0.85
$ws2.Range("D26").Formula = "=B26*C26"

$ws2.Range("A27").Value = "C10uF"
$ws2.Range("B27").Value = 2
$ws2.Range("C27").Value = 0.8
$ws2.Range("D27").Formula = "=B27*C27"

$ws2.Range("A28").Value = "LM1117 3.3V"
$ws2.Range("B28").Value = 2
$ws2.Range("C28").Value = 1.8
$ws2.Range("D28").Formula = "=B28*C28"

$ws2.Range("A29").Value = "LM1117 5V"
$ws2.Range("B29").Value = 2
$ws2.Range("C29").Value = 1.35
$ws2.Range("D29").Formula = "=B29*C29"

# Row30: new "Razem:" (total) row, highlighted yellow
$ws2.Range("A30").Value = "Razem:"
$ws2.Range("D30").Formula = "=B30*C30+SUM(D3:D29)"
$ws2.Range("A30:D30").Interior.Color = 65535

# ---------------------------------------------------------------
# View state: switch active sheet/tab back to Arkusz1, update selections
# ---------------------------------------------------------------
$ws2.Range("D35").Select()
$ws1.Range("D16").Select()
$ws1.Activate()

Write-Output "edit complete"
